$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$tblShape = $s.Shapes.Item(2)
$tbl = $tblShape.Table

# Shift column-2 text up by one row (rows 5-8), and move the freed text into column 1
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "New accounts in last 6 months"
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "No of inquiries"
$tbl.Cell(7, 1).Shape.TextFrame.TextRange.Text = "Primary current balance"
$tbl.Cell(7, 2).Shape.TextFrame.TextRange.Text = "Avg Loan tenure"
$tbl.Cell(8, 1).Shape.TextFrame.TextRange.Text = "Age_at_disbursal"
$tbl.Cell(8, 2).Shape.TextFrame.TextRange.Text = "Time_since_1st_loan"

# Delete the now-redundant trailing rows (originally rows 9 and 10)
$tbl.Rows(10).Delete()
$tbl.Rows(9).Delete()
